$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The notebook/simulation was rerun with two additional samples ("Holden" and
# "Rizzie Spiral") inserted right after "Spiral5" (i.e. as samples #2 and #3).
# This pushes every later sample down by two rows. One of the existing
# samples ("Thomas Hex") was also renamed to "Matthies Hex".
# ---------------------------------------------------------------------------

# Insert two new blank rows for the new samples; this shifts all the old
# rows 4-29 down to rows 6-31, carrying their existing labels/data with them.
$ws.Rows("4:5").Insert()

# Copy the formatting (bold + border) of the column-A index cell above down
# into the two newly inserted column-A cells so they match the rest of the
# table.
$ws.Range("A3").Copy()
$ws.Range("A4:A5").PasteSpecial(-4122)

# Re-number column A (the sample index) for rows 4..31 so it stays sequential
# (2, 3, 4, ..., 29) after the insert.
for ($r = 4; $r -le 31; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 2
}

# New row 4: sample "Holden"
$ws.Cells.Item(4, 2).Value = "Holden"
$ws.Cells.Item(4, 3).Value = 0.8764064587388941
$ws.Cells.Item(4, 4).Value = 1.115620384360406
$ws.Cells.Item(4, 5).Value = 0.9019600965337778
$ws.Cells.Item(4, 6).Value = 0.9019600965337778
$ws.Cells.Item(4, 7).Value = 0.8156850097282374
$ws.Cells.Item(4, 8).Value = 0.8960776325880143
$ws.Cells.Item(4, 9).Value = 0.7330150442468045
$ws.Cells.Item(4, 10).Value = 0.9558039632632692
$ws.Cells.Item(4, 11).Value = 0.9019600965337747
$ws.Cells.Item(4, 12).Value = 2.471025170315269
$ws.Cells.Item(4, 13).Value = 0.7330150442468045
$ws.Cells.Item(4, 14).Value = 0.9019600965337778
$ws.Cells.Item(4, 15).Value = 2.471025170315269
$ws.Cells.Item(4, 16).Value = 1.602020107281037
$ws.Cells.Item(4, 17).Value = 1.673715814527082
$ws.Cells.Item(4, 18).Value = 1.368666770365284
$ws.Cells.Item(4, 19).Value = 1.360148891100323
$ws.Cells.Item(4, 20).Value = 1.368666770365284
$ws.Cells.Item(4, 21).Value = 1.245601692458687
$ws.Cells.Item(4, 22).Value = 1.176873373273705
$ws.Cells.Item(4, 23).Value = 1.095699219971834

# New row 5: sample "Rizzie Spiral"
$ws.Cells.Item(5, 2).Value = "Rizzie Spiral"
$ws.Cells.Item(5, 3).Value = 1.885810824175682
$ws.Cells.Item(5, 4).Value = 1.80426359861646
$ws.Cells.Item(5, 5).Value = 5.643504293258746
$ws.Cells.Item(5, 6).Value = 5.643504293258746
$ws.Cells.Item(5, 7).Value = 1.888167297409928
$ws.Cells.Item(5, 8).Value = 1.886564183279253
$ws.Cells.Item(5, 9).Value = 4.032149763057164
$ws.Cells.Item(5, 10).Value = 0.0008906558968692118
$ws.Cells.Item(5, 11).Value = 5.643504293258746
$ws.Cells.Item(5, 12).Value = 0.7289060913436229
$ws.Cells.Item(5, 13).Value = 4.032149763057164
$ws.Cells.Item(5, 14).Value = 5.643504293258746
$ws.Cells.Item(5, 15).Value = 0.7289060913436229
$ws.Cells.Item(5, 16).Value = 2.380527927200393
$ws.Cells.Item(5, 17).Value = 1.307358457759652
$ws.Cells.Item(5, 18).Value = 3.468186715886511
$ws.Cells.Item(5, 19).Value = 2.215622226192156
$ws.Cells.Item(5, 20).Value = 3.468186715886511
$ws.Cells.Item(5, 21).Value = 3.072592742958804
$ws.Cells.Item(5, 22).Value = 3.586775053018792
$ws.Cells.Item(5, 23).Value = 2.233782088379716

# Rename "Thomas Hex" -> "Matthies Hex" (that sample is now at row 11, since
# it used to be row 9 before the two-row insert shifted it down).
$ws.Cells.Item(11, 2).Value = "Matthies Hex"
